$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.009.64'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.982.56'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '245.71'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.629'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.01'
$ws.Range('E7').Value = '  +3.71%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +1.70%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0800'
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.02'
$ws.Range('E12').Value = '  +9.67%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.845'
$ws.Range('E13').Value = '  +2.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.11'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').Value = '2.277.47'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.47'
$ws.Range('E16').Value = '  +3.94%  '
$ws.Range('D17').Value = '1.989.07'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').Value = '36.886.14'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.19'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.16'
$ws.Range('E21').Value = '  +2.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '230.11'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +1.88%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.152'
$ws.Range('E26').Value = '  +8.65%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.27'
$ws.Range('E27').Value = '  +0.73%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.51'
$ws.Range('E28').Value = '  +2.15%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.55'
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('E30').Value = '  +17.45%  '
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.87'
$ws.Range('E32').Value = '  +2.97%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0621'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.53'
$ws.Range('E34').Value = '  +5.71%  '
$ws.Range('E35').Value = '  +2.72%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.53'
$ws.Range('E39').Value = '  -7.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0980'
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.46'
$ws.Range('E44').Value = '  +3.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '90.25'
$ws.Range('E45').Value = '  +2.90%  '
$ws.Range('D46').Value = '1.368.66'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.25'
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.82'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '46.38'
$ws.Range('E50').Value = '  +6.28%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.96'
$ws.Range('E51').Value = '  +10.30%  '

Write-Host "Done applying changes"
